$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 4498.6  # H113: 4498.3335 -> 4498.6
$ws.Cells.Item(113, 9).Value = 4498.6  # I113: 4498.3335 -> 4498.6
$ws.Cells.Item(113, 11).Value = 4498.6  # K113: 4498.3335 -> 4498.6
$ws.Cells.Item(113, 13).Value = -1244.6  # M113: -1244.3335 -> -1244.6
$ws.Cells.Item(132, 8).Value = 5229.3  # H132: 5340.5 -> 5229.3
$ws.Cells.Item(132, 9).Value = 3731.1667  # I132: 3916.5 -> 3731.1667
$ws.Cells.Item(132, 11).Value = 11193.5001  # K132: 11749.5 -> 11193.5001
$ws.Cells.Item(132, 13).Value = -8663.500100000001  # M132: -9219.5 -> -8663.500100000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3305.5454  # H2: 3446.5833 -> 3305.5454
$ws.Cells.Item(2, 9).Value = 3336.1  # I2: 3487.182 -> 3336.1
$ws.Cells.Item(2, 11).Value = 3336.1  # K2: 3487.182 -> 3336.1
$ws.Cells.Item(2, 13).Value = -3223.1  # M2: -3374.182 -> -3223.1
$ws.Cells.Item(61, 8).Value = 2178  # H61: 2654 -> 2178
$ws.Cells.Item(61, 9).Value = 2501.7  # I61: 3237.125 -> 2501.7
$ws.Cells.Item(61, 11).Value = 2501.7  # K61: 3237.125 -> 2501.7
$ws.Cells.Item(61, 13).Value = -2289.7  # M61: -3025.125 -> -2289.7
$ws.Cells.Item(97, 8).Value = 3296.923  # H97: 3462.5 -> 3296.923
$ws.Cells.Item(97, 9).Value = 2227.5  # I97: 2342.5 -> 2227.5
$ws.Cells.Item(97, 10).Value = 4213.5713  # J97: 4582.5 -> 4213.5713
$ws.Cells.Item(97, 11).Value = 2227.5  # K97: 2342.5 -> 2227.5
$ws.Cells.Item(97, 12).Value = 4213.5713  # L97: 4582.5 -> 4213.5713
$ws.Cells.Item(97, 13).Value = -1731.5  # M97: -1846.5 -> -1731.5
$ws.Cells.Item(97, 14).Value = -5205.5713  # N97: -5574.5 -> -5205.5713
$ws.Cells.Item(116, 8).Value = 3305.5454  # H116: 3446.5833 -> 3305.5454
$ws.Cells.Item(116, 9).Value = 3336.1  # I116: 3487.182 -> 3336.1
$ws.Cells.Item(116, 11).Value = 3336.1  # K116: 3487.182 -> 3336.1
$ws.Cells.Item(116, 13).Value = -1042.1  # M116: -1193.182 -> -1042.1
$ws.Cells.Item(132, 8).Value = 1936.6666  # H132: 2751.5386 -> 1936.6666
$ws.Cells.Item(132, 9).Value = 1932.1428  # I132: 2346.4546 -> 1932.1428
$ws.Cells.Item(132, 10).Value = 2000  # J132: 4979.5 -> 2000
$ws.Cells.Item(132, 11).Value = 5796.428400000001  # K132: 7039.3638 -> 5796.428400000001
$ws.Cells.Item(132, 12).Value = 6000  # L132: 14938.5 -> 6000
$ws.Cells.Item(132, 13).Value = -3266.428400000001  # M132: -4509.3638 -> -3266.428400000001
$ws.Cells.Item(132, 14).Value = -11060  # N132: -19998.5 -> -11060
$ws.Cells.Item(136, 8).Value = 2178  # H136: 2654 -> 2178
$ws.Cells.Item(136, 9).Value = 2501.7  # I136: 3237.125 -> 2501.7
$ws.Cells.Item(136, 11).Value = 7505.099999999999  # K136: 9711.375 -> 7505.099999999999
$ws.Cells.Item(136, 13).Value = -4955.099999999999  # M136: -7161.375 -> -4955.099999999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3305.5454  # H3: 3446.5833 -> 3305.5454
$ws.Cells.Item(3, 9).Value = 3336.1  # I3: 3487.182 -> 3336.1
$ws.Cells.Item(3, 11).Value = 3336.1  # K3: 3487.182 -> 3336.1
$ws.Cells.Item(3, 13).Value = -3222.1  # M3: -3373.182 -> -3222.1
$ws.Cells.Item(20, 8).Value = 2435.875  # H20: 2639.8572 -> 2435.875
$ws.Cells.Item(20, 9).Value = 2715.1428  # I20: 2999.6667 -> 2715.1428
$ws.Cells.Item(20, 11).Value = 2715.1428  # K20: 2999.6667 -> 2715.1428
$ws.Cells.Item(20, 13).Value = -2468.1428  # M20: -2752.6667 -> -2468.1428
$ws.Cells.Item(99, 8).Value = 2800.7693  # H99: 4241.1 -> 2800.7693
$ws.Cells.Item(99, 9).Value = 2951.0833  # I99: 4601.5557 -> 2951.0833
$ws.Cells.Item(99, 11).Value = 2951.0833  # K99: 4601.5557 -> 2951.0833
$ws.Cells.Item(99, 13).Value = -1453.0833  # M99: -3103.5557 -> -1453.0833
$ws.Cells.Item(107, 8).Value = 3778.0588  # H107: 8284.277 -> 3778.0588
$ws.Cells.Item(107, 9).Value = 3441.8  # I107: 3434.5334 -> 3441.8
$ws.Cells.Item(107, 10).Value = 6300  # J107: 32533 -> 6300
$ws.Cells.Item(107, 11).Value = 3441.8  # K107: 3434.5334 -> 3441.8
$ws.Cells.Item(107, 12).Value = 6300  # L107: 32533 -> 6300
$ws.Cells.Item(107, 13).Value = -1521.8  # M107: -1514.5334 -> -1521.8
$ws.Cells.Item(107, 14).Value = -10140  # N107: -36373 -> -10140

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7664.1816  # H31: 8600.888999999999 -> 7664.1816
$ws.Cells.Item(31, 9).Value = 7430.7  # I31: 8426.125 -> 7430.7
$ws.Cells.Item(31, 11).Value = 7430.7  # K31: 8426.125 -> 7430.7
$ws.Cells.Item(31, 13).Value = -7135.7  # M31: -8131.125 -> -7135.7
$ws.Cells.Item(34, 8).Value = 7664.1816  # H34: 8600.888999999999 -> 7664.1816
$ws.Cells.Item(34, 9).Value = 7430.7  # I34: 8426.125 -> 7430.7
$ws.Cells.Item(34, 11).Value = 7430.7  # K34: 8426.125 -> 7430.7
$ws.Cells.Item(34, 13).Value = -7228.7  # M34: -8224.125 -> -7228.7
$ws.Cells.Item(58, 8).Value = 2672.0833  # H58: 2935.0908 -> 2672.0833
$ws.Cells.Item(58, 9).Value = 2808  # I58: 3029.1 -> 2808
$ws.Cells.Item(58, 10).Value = 1992.5  # J58: 1995 -> 1992.5
$ws.Cells.Item(58, 11).Value = 2808  # K58: 3029.1 -> 2808
$ws.Cells.Item(58, 12).Value = 1992.5  # L58: 1995 -> 1992.5
$ws.Cells.Item(58, 13).Value = -2605  # M58: -2826.1 -> -2605
$ws.Cells.Item(58, 14).Value = -2398.5  # N58: -2401 -> -2398.5
$ws.Cells.Item(107, 8).Value = 5481.2  # H107: 5481.7 -> 5481.2
$ws.Cells.Item(107, 9).Value = 2468.8333  # I107: 2469.6667 -> 2468.8333
$ws.Cells.Item(107, 11).Value = 2468.8333  # K107: 2469.6667 -> 2468.8333
$ws.Cells.Item(107, 13).Value = -548.8332999999998  # M107: -549.6667000000002 -> -548.8332999999998
$ws.Cells.Item(122, 8).Value = 1027.1  # H122: 1412.625 -> 1027.1
$ws.Cells.Item(122, 9).Value = 1035.25  # I122: 1412.625 -> 1035.25
$ws.Cells.Item(122, 10).Value = 994.5  # J122: 0 -> 994.5
$ws.Cells.Item(122, 11).Value = 3105.75  # K122: 4237.875 -> 3105.75
$ws.Cells.Item(122, 12).Value = 2983.5  # L122: 0 -> 2983.5
$ws.Cells.Item(122, 13).Value = -655.75  # M122: -1787.875 -> -655.75
$ws.Cells.Item(122, 14).Value = -7883.5  # N122: None -> -7883.5
$ws.Cells.Item(132, 8).Value = 1743.7  # H132: 1893.1111 -> 1743.7
$ws.Cells.Item(132, 9).Value = 1137.4  # I132: 1322 -> 1137.4
$ws.Cells.Item(132, 11).Value = 3412.2  # K132: 3966 -> 3412.2
$ws.Cells.Item(132, 13).Value = -882.2000000000003  # M132: -1436 -> -882.2000000000003
$ws.Cells.Item(134, 8).Value = 2525.4614  # H134: 2566.7273 -> 2525.4614
$ws.Cells.Item(134, 9).Value = 2484.4  # I134: 2530.875 -> 2484.4
$ws.Cells.Item(134, 11).Value = 7453.200000000001  # K134: 7592.625 -> 7453.200000000001
$ws.Cells.Item(134, 13).Value = -4918.200000000001  # M134: -5057.625 -> -4918.200000000001
$ws.Cells.Item(136, 8).Value = 2672.0833  # H136: 2935.0908 -> 2672.0833
$ws.Cells.Item(136, 9).Value = 2808  # I136: 3029.1 -> 2808
$ws.Cells.Item(136, 10).Value = 1992.5  # J136: 1995 -> 1992.5
$ws.Cells.Item(136, 11).Value = 8424  # K136: 9087.299999999999 -> 8424
$ws.Cells.Item(136, 12).Value = 5977.5  # L136: 5985 -> 5977.5
$ws.Cells.Item(136, 13).Value = -5874  # M136: -6537.299999999999 -> -5874
$ws.Cells.Item(136, 14).Value = -11077.5  # N136: -11085 -> -11077.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 570  # H11: 500 -> 570
$ws.Cells.Item(11, 9).Value = 570  # I11: 500 -> 570
$ws.Cells.Item(11, 11).Value = 1710  # K11: 1500 -> 1710
$ws.Cells.Item(11, 13).Value = -1570  # M11: -1360 -> -1570

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5000  # H70: 4993.5 -> 5000
$ws.Cells.Item(70, 9).Value = 0  # I70: 4987 -> 0
$ws.Cells.Item(70, 11).Value = 0  # K70: 4987 -> 0
$ws.Cells.Item(70, 13).ClearContents()  # M70: -4717 -> (removed)
$ws.Cells.Item(73, 8).Value = 5000  # H73: 4993.5 -> 5000
$ws.Cells.Item(73, 9).Value = 0  # I73: 4987 -> 0
$ws.Cells.Item(73, 11).Value = 0  # K73: 4987 -> 0
$ws.Cells.Item(73, 13).ClearContents()  # M73: -4051 -> (removed)
$ws.Cells.Item(102, 8).Value = 3166.8  # H102: 3240.889 -> 3166.8
$ws.Cells.Item(102, 10).Value = 3599  # J102: 4698 -> 3599
$ws.Cells.Item(102, 12).Value = 3599  # L102: 4698 -> 3599
$ws.Cells.Item(102, 14).Value = -6843  # N102: -7942 -> -6843
$ws.Cells.Item(113, 8).Value = 816.6667  # H113: 875 -> 816.6667
$ws.Cells.Item(113, 9).Value = 816.6667  # I113: 875 -> 816.6667
$ws.Cells.Item(113, 11).Value = 816.6667  # K113: 875 -> 816.6667
$ws.Cells.Item(113, 13).Value = 1353.3333  # M113: 1295 -> 1353.3333
$ws.Cells.Item(122, 8).Value = 4624.5  # H122: 6201.6 -> 4624.5
$ws.Cells.Item(122, 9).Value = 5997.25  # I122: 10000 -> 5997.25
$ws.Cells.Item(122, 10).Value = 3251.75  # J122: 3669.3333 -> 3251.75
$ws.Cells.Item(122, 11).Value = 17991.75  # K122: 30000 -> 17991.75
$ws.Cells.Item(122, 12).Value = 9755.25  # L122: 11007.9999 -> 9755.25
$ws.Cells.Item(122, 13).Value = -15541.75  # M122: -27550 -> -15541.75
$ws.Cells.Item(122, 14).Value = -14655.25  # N122: -15907.9999 -> -14655.25
$ws.Cells.Item(126, 8).Value = 6310.5557  # H126: 6816.1665 -> 6310.5557
$ws.Cells.Item(126, 9).Value = 5849.5  # I126: 6179.6 -> 5849.5
$ws.Cells.Item(126, 11).Value = 17548.5  # K126: 18538.8 -> 17548.5
$ws.Cells.Item(126, 13).Value = -15078.5  # M126: -16068.8 -> -15078.5
$ws.Cells.Item(132, 8).Value = 3106.9412  # H132: 3238.625 -> 3106.9412
$ws.Cells.Item(132, 9).Value = 2721.3333  # I132: 2844.2856 -> 2721.3333
$ws.Cells.Item(132, 11).Value = 8163.999899999999  # K132: 8532.856800000001 -> 8163.999899999999
$ws.Cells.Item(132, 13).Value = -5633.999899999999  # M132: -6002.856800000001 -> -5633.999899999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4164.3076  # H40: 4684.4546 -> 4164.3076
$ws.Cells.Item(40, 9).Value = 3792.889  # I40: 4504.143 -> 3792.889
$ws.Cells.Item(40, 11).Value = 3792.889  # K40: 4504.143 -> 3792.889
$ws.Cells.Item(40, 13).Value = -3656.889  # M40: -4368.143 -> -3656.889
$ws.Cells.Item(122, 8).Value = 5473.5  # H122: 4372.5835 -> 5473.5
$ws.Cells.Item(122, 9).Value = 4557.6  # I122: 3496.7778 -> 4557.6
$ws.Cells.Item(122, 11).Value = 13672.8  # K122: 10490.3334 -> 13672.8
$ws.Cells.Item(122, 13).Value = -11222.8  # M122: -8040.3334 -> -11222.8
$ws.Cells.Item(132, 8).Value = 2046  # H132: 1871.6666 -> 2046
$ws.Cells.Item(132, 10).Value = 0  # J132: 1000 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 3000 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -8060 -> (removed)
$ws.Cells.Item(136, 8).Value = 3489.6365  # H136: 3579.2727 -> 3489.6365
$ws.Cells.Item(136, 9).Value = 4410.75  # I136: 4534 -> 4410.75
$ws.Cells.Item(136, 11).Value = 13232.25  # K136: 13602 -> 13232.25
$ws.Cells.Item(136, 13).Value = -10682.25  # M136: -11052 -> -10682.25

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2049.0667  # H126: 2526.3333 -> 2049.0667
$ws.Cells.Item(126, 9).Value = 1653.1  # I126: 2088.5 -> 1653.1
$ws.Cells.Item(126, 10).Value = 2841  # J126: 3402 -> 2841
$ws.Cells.Item(126, 11).Value = 4959.299999999999  # K126: 6265.5 -> 4959.299999999999
$ws.Cells.Item(126, 12).Value = 8523  # L126: 10206 -> 8523
$ws.Cells.Item(126, 13).Value = -2489.299999999999  # M126: -3795.5 -> -2489.299999999999
$ws.Cells.Item(126, 14).Value = -13463  # N126: -15146 -> -13463
$ws.Cells.Item(132, 8).Value = 3372.5386  # H132: 3529.5 -> 3372.5386
$ws.Cells.Item(132, 9).Value = 2242.2632  # I132: 2284.111 -> 2242.2632
$ws.Cells.Item(132, 10).Value = 6440.4287  # J132: 7265.6665 -> 6440.4287
$ws.Cells.Item(132, 11).Value = 6726.7896  # K132: 6852.333 -> 6726.7896
$ws.Cells.Item(132, 12).Value = 19321.2861  # L132: 21796.9995 -> 19321.2861
$ws.Cells.Item(132, 13).Value = -4196.7896  # M132: -4322.333 -> -4196.7896
$ws.Cells.Item(132, 14).Value = -24381.2861  # N132: -26856.9995 -> -24381.2861
$ws.Cells.Item(136, 8).Value = 1403.0769  # H136: 1438.52 -> 1403.0769
$ws.Cells.Item(136, 9).Value = 1436.4584  # I136: 1476.4348 -> 1436.4584
$ws.Cells.Item(136, 11).Value = 4309.3752  # K136: 4429.3044 -> 4309.3752
$ws.Cells.Item(136, 13).Value = -1759.3752  # M136: -1879.3044 -> -1759.3752
